$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 19 data
$ws.Range("A19").Value = 45759
$ws.Range("B19").Value = 74
$ws.Range("C19").Value = 75
$ws.Range("D19").Value = 70

# Number formats: row 18 (no longer last) gets plain date format,
# row 19 (new last row) gets the "latest" date-time format
$ws.Range("A18").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("A19").NumberFormat = "YYYY-MM-DD"
